# Slide 7 ("HTTPS, HTTP2"): highlight part of the screenshot with a red
# outlined rectangle, then group that rectangle together with the picture.
#
# The highlight rectangle is copy/pasted from the equivalent "직사각형 2"
# shape that already exists on slide 5 of this deck (same red-outline /
# no-fill / shape-style / empty-centered-text template the author reuses
# throughout the deck) so the resulting <p:style>/<p:txBody> markup matches
# exactly, then it is repositioned/resized and grouped with the picture.

$p = $ppt.ActivePresentation
$srcSlide = $p.Slides.Item(5)
$slide = $p.Slides.Item(7)

$pic = $slide.Shapes.Item(1)

# Copy the template highlight-rectangle shape from slide 5 and paste it
# onto slide 7.
$template = $srcSlide.Shapes.Item(2)
$template.Copy()
$pasted = $slide.Shapes.Paste()
$rect = $pasted.Item(1)

# Position/size the rectangle over the "HTTPS, HTTP2" area of the picture
# (values are EMU/12700, chosen so the stored EMU survives the COM
# single-precision round-trip exactly).
$rect.Left = 367.78377952755903
$rect.Top = 402.8107874015748
$rect.Width = 194.59456692913386
$rect.Height = 51.567560255118124

# Group the picture and the new rectangle together.
$range = $slide.Shapes.Range(@(1, 2))
$grp = $range.Group()
$grp.Name = "그룹 3"

# Reposition the resulting group (its children keep their original
# coordinates; only the group's own offset moves).
$grp.Left = 191.7315751031496
$grp.Top = 48.25504117007876

Write-Output "Grouped picture + highlight rectangle on slide 7 ($($grp.Name), $($grp.GroupItems.Count) items)"
